# Generate Report for Handoff
#
# Status moves from "In Translation" to "Ready for handoff" and the
# "Latest Handoff" / "Latest HO Xliff Generate Date" timestamps advance to
# reflect the newly generated handoff report. The affected status/date
# columns are also widened slightly to fit the new "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-12 23:14:47"

# Widen the zh-cn / de-de status columns to fit "Ready for handoff"
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-12 23:14:40"
$zhcn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-12 23:14:47"
$dede.Columns.Item(3).ColumnWidth = 16.3
